# The commit swaps the colour palette that used to live in ppt/theme/theme1.xml
# ("Integral") for the stock Office palette that used to live in
# ppt/theme/theme2.xml ("Office Theme") -- the font scheme and format scheme
# are already byte-identical between the two theme parts, so only the 12
# colour-scheme slots (clrScheme) actually change for the part that is
# reachable through the PowerPoint object model (the deck's primary theme,
# i.e. the slide master's theme / theme1.xml).
#
# msoThemeColor slot order (ThemeColorScheme.Colors / .Item index 1-12):
#   1 dk1   2 lt1   3 dk2   4 lt2
#   5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#   11 hlink   12 folHlink

function ConvertTo-BGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office Theme" colour scheme (RRGGBB, matches the a:clrScheme the
# diff introduces).
$officeTheme = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = ConvertTo-BGR $officeTheme[$i]
}

Write-Output "Theme colour scheme updated to Office Theme palette."
